# Add active tool functionality
#
# - tool_checkout_log: append new sign-out/sign-in rows (rows 2-12)
# - tools: add a "Status" column, marking currently-checked-out tools "Active"
# - selection / active-sheet bookkeeping to match the final UI state

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. tool_checkout_log — rewrite the log with the fuller history
# ---------------------------------------------------------------------------
$log = $wb.Worksheets.Item("tool_checkout_log")

$logData = @(
    @("emp1", "tool1", "02/10/2024 14:25", "02/10/2024 14:27", "emp1"),
    @("emp2", "tool2", "02/10/2024 14:26", "", ""),
    @("emp3", "tool3", "02/10/2024 14:26", "02/10/2024 14:40", "emp3"),
    @("emp4", "tool4", "02/10/2024 14:26", "02/10/2024 14:27", "emp4"),
    @("emp5", "tool5", "02/10/2024 14:26", "02/10/2024 14:37", "emp4"),
    @("emp6", "tool6", "02/10/2024 14:26", "02/10/2024 14:39", "emp6"),
    @("emp7", "tool7", "02/10/2024 14:26", "02/10/2024 14:27", "emp7"),
    @("emp8", "tool8", "02/10/2024 14:26", "02/10/2024 14:27", "emp8"),
    @("emp9", "tool9", "02/10/2024 14:26", "", ""),
    @("emp5", "tool5", "02/10/2024 14:39", "02/10/2024 14:43", "emp5"),
    @("emp7", "tool7", "02/10/2024 14:40", "", "")
)

for ($i = 0; $i -lt $logData.Count; $i++) {
    $r = $i + 2
    $row = $logData[$i]
    $log.Cells.Item($r, 1).Value = $row[0]
    $log.Cells.Item($r, 2).Value = $row[1]
    $log.Cells.Item($r, 3).Value = $row[2]
    if ($row[3] -ne "") {
        $log.Cells.Item($r, 4).Value = $row[3]
    } else {
        $log.Cells.Item($r, 4).ClearContents()
    }
    if ($row[4] -ne "") {
        $log.Cells.Item($r, 5).Value = $row[4]
    } else {
        $log.Cells.Item($r, 5).ClearContents()
    }
}

# narrow the Sign Out Time / Sign In Time columns now that every row is filled in,
# and let the Sign Out Employee / Sign In Employee columns fall back to the sheet's
# default (un-customized) width
$log.Columns.Item(1).ColumnWidth = 8
$log.Columns.Item(3).ColumnWidth = 16
$log.Columns.Item(4).ColumnWidth = 15.83
$log.Columns.Item(5).ColumnWidth = 8

# ---------------------------------------------------------------------------
# 2. tools — add the "Status" column
# ---------------------------------------------------------------------------
$tools = $wb.Worksheets.Item("tools")
$tools.Range("C1").Value = "Status"
$tools.Range("C3").Value = "Active"
$tools.Range("C8").Value = "Active"
$tools.Range("C10").Value = "Active"

# ---------------------------------------------------------------------------
# 3. selection bookkeeping to match the saved UI state
# ---------------------------------------------------------------------------
$sheet1 = $wb.Worksheets.Item("Sheet")
[void]$sheet1.Range("A1").Select()

[void]$log.Range("A1").Select()

$employees = $wb.Worksheets.Item("employees")
$employees.Columns.Item(1).ColumnWidth = 8
[void]$employees.Activate()
[void]$employees.Range("B2:B11").Select()

[void]$tools.Activate()
[void]$tools.Range("H40").Select()
